$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column D holds text-formatted numbers (e.g. "26.429.38" with
# thousands separators). Force the whole column to Text format first so
# Excel doesn't silently re-interpret plain-looking numeric strings
# (like "306.56") as actual numbers when we assign them below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.429.38"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.810.60"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "306.56"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "0.3585"
$ws.Range("E8").Value = "  -1.66%  "
$ws.Range("D9").Value = "46.48"
$ws.Range("E9").Value = "  +4.04%  "
$ws.Range("D10").Value = "0.07074"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").Value = "0.8938"
$ws.Range("E11").Value = "  +3.19%  "
$ws.Range("D12").Value = "0.07808"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "19.37"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").Value = "1.799.72"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").Value = "5.286"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "6.307"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "85.10"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D19").Value = "0.000008512"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "26.471.78"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "4.970"
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("D24").Value = "2.040.03"
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("D25").Value = "10.51"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").Value = "1.961"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").Value = "152.01"
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("D28").Value = "17.80"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Value = "2.054"
$ws.Range("E29").Value = "  +4.33%  "
$ws.Range("D30").Value = "112.37"
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "0.08695"
$ws.Range("E32").Value = "  +0.89%  "
$ws.Range("D33").Value = "3.108"
$ws.Range("E33").Value = "  +2.77%  "
$ws.Range("D34").Value = "2.792"
$ws.Range("E34").Value = "  +10.14%  "
$ws.Range("D35").Value = "4.456"
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("D36").Value = "0.7273"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "1.076"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").Value = "0.01927"
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("D40").Value = "0.05119"
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("D41").Value = "2.897"
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("D42").Value = "0.5085"
$ws.Range("E42").Value = "  +3.95%  "
$ws.Range("D43").Value = "6.772"
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("D44").Value = "0.1515"
$ws.Range("E44").Value = "  -3.02%  "
$ws.Range("D45").Value = "8.039"
$ws.Range("E45").Value = "  -0.65%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  +1.41%  "
$ws.Range("D48").Value = "10.02"
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("D49").Value = "99.86"
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("D50").Value = "1.572"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").Value = "0.05988"
$ws.Range("E51").Value = "  +0.09%  "